$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.437.92'
$ws.Range('E2').Value = '  +0.56%  '
$ws.Range('D3').Value = '2.238.47'
$ws.Range('E3').Value = '  -0.49%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.24'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E6').Value = '  -0.77%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '74.27'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.77%  '
$ws.Range('E8').Value = '  +0.16%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.622'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.53%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '43.32'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.00%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0961'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.86%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.12'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.10%  '
$ws.Range('E13').Value = '  -0.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.47'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.54%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.853'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.02%  '
$ws.Range('D16').Value = '2.208.47'
$ws.Range('E16').Value = '  -1.02%  '
$ws.Range('D17').Value = '42.256.02'
$ws.Range('E17').Value = '  +0.36%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0000108'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +9.77%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.17'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.84%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.04'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.38'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +39.05%  '
$ws.Range('E22').Value = '  -0.18%  '
$ws.Range('E23').Value = '  -5.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.80'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.81%  '
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('E26').Value = '  +0.96%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.30'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('E28').Value = '  +7.18%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '166.74'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.86%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '20.96'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.51%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.84'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +18.38%  '
$ws.Range('E32').Value = '  -3.23%  '
$ws.Range('E33').Value = '  -2.34%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '30.00'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -8.94%  '
$ws.Range('E35').Value = '  -0.41%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.47'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.32%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0310'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.82%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '13.23'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -7.34%  '
$ws.Range('E39').Value = '  -0.96%  '
$ws.Range('E40').Value = '  -4.46%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '63.64'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.12%  '
$ws.Range('E42').Value = '  -0.80%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.93'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.82%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '106.09'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -6.33%  '
$ws.Range('E45').Value = '  +2.65%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.995'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.26%  '
$ws.Range('E47').Value = '  +0.92%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.39'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.63%  '
$ws.Range('E49').Value = '  +0.68%  '
$ws.Range('E50').Value = '  +1.39%  '
$ws.Range('E51').Value = '  -2.52%  '
